# edit.ps1 -- applies the "about.docx" change set described in the task diff:
#   1. Split the two "<space><e-mail>" hyperlink runs into a leading-space
#      run and a separate e-mail run (same Hyperlink character style).
#   2. Insert a new "Abstract Title" paragraph style (next-style: Abstract).
#   3. Change the "Abstract" style's space-before from 15pt (300) to 5pt (100).
#   4. Insert a new "Footnote Block Text" paragraph style, based on
#      "Footnote Text".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1a. "AMČR Praha" contact block: " amcr@arup.cas.cz" -> " " + "amcr@arup.cas.cz"
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("amcr@arup.cas.cz", $true, $false, $false, $false, `
                           $false, $true, 1, $false, $null, 0)
# Toggle a character attribute on and back off on just the e-mail text (not
# the leading space) -- this forces the run to split from its neighbour
# while leaving the resulting formatting identical to the original.
$rng.Font.Bold = $true
$rng.Font.Bold = $false

# ---------------------------------------------------------------------------
# 1b. "AMČR Brno" contact block: " amcr@arub.cz" -> " " + "amcr@arub.cz"
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$null = $rng2.Find.Execute("amcr@arub.cz", $true, $false, $false, $false, `
                            $false, $true, 1, $false, $null, 0)
$rng2.Font.Bold = $true
$rng2.Font.Bold = $false

# ---------------------------------------------------------------------------
# 2. New paragraph style "Abstract Title" (mirrors pandoc's reference.docx)
# ---------------------------------------------------------------------------
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = $d.Styles.Item("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles.Item("Abstract")
$abstractTitle.QuickStyle = $true
$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060

# ---------------------------------------------------------------------------
# 3. "Abstract" style: space-before 300 (15pt) -> 100 (5pt)
# ---------------------------------------------------------------------------
$abstract = $d.Styles.Item("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# ---------------------------------------------------------------------------
# 4. New paragraph style "Footnote Block Text", based on "Footnote Text"
# ---------------------------------------------------------------------------
$footnoteBlockText = $d.Styles.Add("Footnote Block Text", 1)
$footnoteBlockText.BaseStyle = $d.Styles.Item("Footnote Text")
$footnoteBlockText.NextParagraphStyle = $d.Styles.Item("Footnote Text")
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true
$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24

Write-Output "edit.ps1 completed"
